# Auto-generated: apply 2023-10-10 data updates across all affected worksheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 5946
$ws.Range("J3").Value = 6365
$ws.Range("G4").Value = 1473
$ws.Range("J4").Value = 1375
$ws.Range("J5").Value = 488
$ws.Range("J6").Value = 8200
$ws.Range("G7").Value = 24698
$ws.Range("J7").Value = 22374

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 383
$ws.Range("J3").Value = 428
$ws.Range("J6").Value = 482
$ws.Range("J7").Value = 1406

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 244
$ws.Range("G4").Value = 49
$ws.Range("J4").Value = 44
$ws.Range("J6").Value = 359
$ws.Range("G7").Value = 1426
$ws.Range("J7").Value = 1033

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J3").Value = 116
$ws.Range("J6").Value = 73
$ws.Range("J7").Value = 327

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J6").Value = 198
$ws.Range("J7").Value = 685

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J5").Value = 16
$ws.Range("J6").Value = 198
$ws.Range("J7").Value = 559

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J6").Value = 163
$ws.Range("J7").Value = 663
$ws.Range("J8").Value = 1406
$ws.Range("J10").Value = 158
$ws.Range("J11").Value = 356
$ws.Range("J15").Value = 247
$ws.Range("J16").Value = 89
$ws.Range("J19").Value = 664
$ws.Range("J20").Value = 469
$ws.Range("J25").Value = 111
$ws.Range("J29").Value = 1234
$ws.Range("J31").Value = 201
$ws.Range("G33").Value = 1426
$ws.Range("J33").Value = 1033
$ws.Range("J36").Value = 309
$ws.Range("J37").Value = 685
$ws.Range("J41").Value = 148
$ws.Range("J42").Value = 946
$ws.Range("J43").Value = 181
$ws.Range("J47").Value = 170
$ws.Range("J51").Value = 280
$ws.Range("J52").Value = 562
$ws.Range("J54").Value = 434
$ws.Range("J55").Value = 309
$ws.Range("J63").Value = 78
$ws.Range("J65").Value = 559
$ws.Range("J67").Value = 854
$ws.Range("J73").Value = 215
$ws.Range("J76").Value = 341
$ws.Range("J79").Value = 639
$ws.Range("J84").Value = 189
$ws.Range("J85").Value = 919
$ws.Range("J89").Value = 293
$ws.Range("J90").Value = 240
$ws.Range("J91").Value = 255
$ws.Range("J92").Value = 70
$ws.Range("J93").Value = 98
$ws.Range("J94").Value = 231
$ws.Range("J95").Value = 327
$ws.Range("J96").Value = 256
$ws.Range("J97").Value = 194
$ws.Range("J98").Value = 162
$ws.Range("G101").Value = 24698
$ws.Range("J101").Value = 22374

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J6").Value = 56
$ws.Range("J7").Value = 201

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 212
$ws.Range("J3").Value = 323
$ws.Range("J4").Value = 64
$ws.Range("J6").Value = 230
$ws.Range("J7").Value = 854

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J6").Value = 58
$ws.Range("J7").Value = 189

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J6").Value = 209
$ws.Range("J7").Value = 434

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 378
$ws.Range("J5").Value = 46
$ws.Range("J6").Value = 313
$ws.Range("J7").Value = 1234

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J3").Value = 195
$ws.Range("J6").Value = 255
$ws.Range("J7").Value = 664

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J2").Value = 54
$ws.Range("J7").Value = 341

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J2").Value = 47
$ws.Range("J7").Value = 163

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J6").Value = 85
$ws.Range("J7").Value = 148

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 203
$ws.Range("J6").Value = 493
$ws.Range("J7").Value = 946

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J6").Value = 85
$ws.Range("J7").Value = 158

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J2").Value = 67
$ws.Range("J4").Value = 13
$ws.Range("J6").Value = 160
$ws.Range("J7").Value = 309

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J3").Value = 70
$ws.Range("J6").Value = 90
$ws.Range("J7").Value = 256

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J3").Value = 107
$ws.Range("J7").Value = 255

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 180
$ws.Range("J3").Value = 219
$ws.Range("J6").Value = 185
$ws.Range("J7").Value = 639

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J2").Value = 130
$ws.Range("J6").Value = 126
$ws.Range("J7").Value = 469

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J6").Value = 95
$ws.Range("J7").Value = 309

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("J6").Value = 35
$ws.Range("J7").Value = 98

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 202
$ws.Range("J6").Value = 215
$ws.Range("J7").Value = 663

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J2").Value = 42
$ws.Range("J6").Value = 128
$ws.Range("J7").Value = 231

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J2").Value = 49
$ws.Range("J7").Value = 111

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J6").Value = 80
$ws.Range("J7").Value = 170

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J3").Value = 58
$ws.Range("J6").Value = 104
$ws.Range("J7").Value = 247

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J3").Value = 25
$ws.Range("J7").Value = 162

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J2").Value = 107
$ws.Range("J6").Value = 151
$ws.Range("J7").Value = 356

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J3").Value = 57
$ws.Range("J7").Value = 215

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J6").Value = 135
$ws.Range("J7").Value = 194

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("J2").Value = 20
$ws.Range("J7").Value = 70

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J6").Value = 89
$ws.Range("J7").Value = 293

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J5").Value = 9
$ws.Range("J7").Value = 240

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J6").Value = 110
$ws.Range("J7").Value = 280

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J6").Value = 106
$ws.Range("J7").Value = 181

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 244
$ws.Range("J3").Value = 329
$ws.Range("J7").Value = 919

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J3").Value = 171
$ws.Range("J6").Value = 228
$ws.Range("J7").Value = 562

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("J2").Value = 11
$ws.Range("J7").Value = 89
